# Fruta / hortaliza, semanal
# Insert a new weekly price block (date 2022-10-25, serial 44859) right after the
# existing top block (rows 2-4), pushing all subsequent blocks down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 5 downward by 3 rows, preserving formatting/styles.
$ws.Rows("5:7").Insert()

# New data block for 2022-10-25 (serial 44859) - Frutilla, Agrícola del Norte S.A. de Arica
$newBlock = @(
    @(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44859, 15, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Primera", 550, 6000, 7000, 6545, "`$/bandeja 3 kilos", "Región de Arica y Parinacota", 2182, 3),
    @(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44859, 15, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Segunda", 500, 5000, 6000, 5600, "`$/bandeja 3 kilos", "Región de Arica y Parinacota", 1867, 3),
    @(1, "Agrícola del Norte S.A. de Arica", "Arica y Parinacota", 44859, 15, "Fruta", 100101, "Berries", 100112025, "Frutilla", "Sin especificar", "Tercera", 350, 4000, 5000, 4857, "`$/bandeja 3 kilos", "Región de Arica y Parinacota", 1619, 3)
)

$startRow = 5
for ($i = 0; $i -lt $newBlock.Length; $i++) {
    $rowValues = $newBlock[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
